# Auto-generated edit script: updates cached market-price / profit
# figures across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables
# to reflect a refreshed data pull from the market board.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2462.1177
$ws.Range("I17").Value = 2875
$ws.Range("J17").Value = 1997.625
$ws.Range("K17").Value = 8625
$ws.Range("L17").Value = 5992.875
$ws.Range("M17").Value = -8457
$ws.Range("N17").Value = -6328.875
$ws.Range("H138").Value = 3950.7869
$ws.Range("I138").Value = 1218.6666
$ws.Range("J138").Value = 4248.8364
$ws.Range("K138").Value = 3655.9998
$ws.Range("L138").Value = 12746.5092
$ws.Range("M138").Value = 1484.0002
$ws.Range("N138").Value = -23026.5092

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 396.66666
$ws.Range("I2").Value = 389.2857
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 389.2857
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -276.2857
$ws.Range("N2").Value = -726
$ws.Range("H52").Value = 10709
$ws.Range("I52").Value = 10709
$ws.Range("K52").Value = 10709
$ws.Range("M52").Value = -10391
$ws.Range("H61").Value = 1809
$ws.Range("I61").Value = 1712.2222
$ws.Range("J61").Value = 2244.5
$ws.Range("K61").Value = 1712.2222
$ws.Range("L61").Value = 2244.5
$ws.Range("M61").Value = -1500.2222
$ws.Range("N61").Value = -2668.5
$ws.Range("H63").Value = 4946.7617
$ws.Range("I63").Value = 3903.0908
$ws.Range("K63").Value = 3903.0908
$ws.Range("M63").Value = -3217.0908
$ws.Range("H66").Value = 4946.7617
$ws.Range("I66").Value = 3903.0908
$ws.Range("K66").Value = 19515.454
$ws.Range("M66").Value = -16083.454
$ws.Range("H102").Value = 1242.8572
$ws.Range("I102").Value = 640
$ws.Range("K102").Value = 640
$ws.Range("M102").Value = 982
$ws.Range("H116").Value = 396.66666
$ws.Range("I116").Value = 389.2857
$ws.Range("J116").Value = 500
$ws.Range("K116").Value = 389.2857
$ws.Range("L116").Value = 500
$ws.Range("M116").Value = 1904.7143
$ws.Range("N116").Value = -5088
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 2182.5
$ws.Range("I132").Value = 1952.3334
$ws.Range("K132").Value = 5857.0002
$ws.Range("M132").Value = -3327.0002
$ws.Range("H136").Value = 1809
$ws.Range("I136").Value = 1712.2222
$ws.Range("J136").Value = 2244.5
$ws.Range("K136").Value = 5136.6666
$ws.Range("L136").Value = 6733.5
$ws.Range("M136").Value = -2586.6666
$ws.Range("N136").Value = -11833.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 396.66666
$ws.Range("I3").Value = 389.2857
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 389.2857
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = -275.2857
$ws.Range("N3").Value = -728
$ws.Range("H50").Value = 27500
$ws.Range("J50").Value = 27500
$ws.Range("L50").Value = 27500
$ws.Range("N50").Value = -28648
$ws.Range("H105").Value = 3685.2593
$ws.Range("I105").Value = 3012.1904
$ws.Range("K105").Value = 3012.1904
$ws.Range("M105").Value = -1265.1904
$ws.Range("H107").Value = 824.4666999999999
$ws.Range("I107").Value = 824.4666999999999
$ws.Range("K107").Value = 824.4666999999999
$ws.Range("M107").Value = 1095.5333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7041.143
$ws.Range("I31").Value = 3190.5
$ws.Range("J31").Value = 8581.4
$ws.Range("K31").Value = 3190.5
$ws.Range("L31").Value = 8581.4
$ws.Range("M31").Value = -2895.5
$ws.Range("N31").Value = -9171.4
$ws.Range("H34").Value = 7041.143
$ws.Range("I34").Value = 3190.5
$ws.Range("J34").Value = 8581.4
$ws.Range("K34").Value = 3190.5
$ws.Range("L34").Value = 8581.4
$ws.Range("M34").Value = -2988.5
$ws.Range("N34").Value = -8985.4
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 5000
$ws.Range("M76").Value = -4685
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 5000
$ws.Range("M79").Value = -3908
$ws.Range("H132").Value = 4049.6924
$ws.Range("I132").Value = 2654
$ws.Range("J132").Value = 5246
$ws.Range("K132").Value = 7962
$ws.Range("L132").Value = 15738
$ws.Range("M132").Value = -5432
$ws.Range("N132").Value = -20798

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 65210330
$ws.Range("I4").Value = 84229360
$ws.Range("K4").Value = 252688080
$ws.Range("M4").Value = -252687968
$ws.Range("H34").Value = 1243.8462
$ws.Range("I34").Value = 730.3333
$ws.Range("J34").Value = 1684
$ws.Range("K34").Value = 2190.9999
$ws.Range("L34").Value = 5052
$ws.Range("M34").Value = -2106.9999
$ws.Range("N34").Value = -5220
$ws.Range("H132").Value = 1323.75
$ws.Range("J132").Value = 1269.75
$ws.Range("L132").Value = 11427.75
$ws.Range("N132").Value = -16487.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3449.875
$ws.Range("I132").Value = 2465.25
$ws.Range("K132").Value = 7395.75
$ws.Range("M132").Value = -4865.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3976.5898
$ws.Range("I132").Value = 3636
$ws.Range("J132").Value = 4466.1875
$ws.Range("K132").Value = 10908
$ws.Range("L132").Value = 13398.5625
$ws.Range("M132").Value = -8378
$ws.Range("N132").Value = -18458.5625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 34990
$ws.Range("I51").Value = 34990
$ws.Range("K51").Value = 34990
$ws.Range("M51").Value = -34480
$ws.Range("H81").Value = 1115.375
$ws.Range("I81").Value = 1115.375
$ws.Range("K81").Value = 2230.75
$ws.Range("M81").Value = -1169.75
$ws.Range("H84").Value = 1115.375
$ws.Range("I84").Value = 1115.375
$ws.Range("K84").Value = 11153.75
$ws.Range("M84").Value = -5849.75
$ws.Range("H100").Value = 3168.4285
$ws.Range("I100").Value = 3295.8
$ws.Range("J100").Value = 2850
$ws.Range("K100").Value = 6591.6
$ws.Range("L100").Value = 5700
$ws.Range("M100").Value = -6050.6
$ws.Range("N100").Value = -6782
$ws.Range("H132").Value = 1548.1052
$ws.Range("I132").Value = 1384.1111
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 4152.3333
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -1622.3333
$ws.Range("N132").Value = -18560
